$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# Insert a new row at row 4, pushing existing rows 4-11 down to 5-12
$ws.Rows.Item(4).Insert()

# Fill the new row 4 with the climate change factor variable
$ws.Cells.Item(4, 1).Value = "General"
$ws.Cells.Item(4, 2).Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = ""
$ws.Cells.Item(4, 7).Value = ""
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 0.5

for ($col = 10; $col -le 45; $col++) {
    $ws.Cells.Item(4, $col).Value = 1
}
